# "suppression page etude complementaire" - remove the "Operation/constraint"
# row contents on the Include sheet and replace them with "Codes / All codes",
# point the System URI at EDQM Standard Terms instead of SNOMED CT, and bump
# the Date metadata value on the Metadata sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": refresh generation Date ---
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Range("B8").Value = "2025-04-09T16:24:06+00:00"

# --- Sheet "Include #0": drop the filter/constraint columns, switch to "all codes" ---
$wsInclude = $wb.Worksheets.Item(2)

# Row 1/2, column A: "Property"/"constraint" -> "Codes"/"All codes"
$wsInclude.Range("A1").Value = "Codes"
$wsInclude.Range("A2").Value = "All codes"

# Row 1/2, column B ("Operation"/"=") and column C ("Value"/the SNOMED
# expression) are no longer used at all - clear them completely (including
# formatting) so the cells disappear from the sheet and the used range
# shrinks from A1:C4 down to A1:B4.
$wsInclude.Range("B1").Clear()
$wsInclude.Range("B2").Clear()
$wsInclude.Range("C1:C2").Clear()

# Row 4: keep the "System URI" label, but point it at EDQM Standard Terms
# instead of SNOMED CT.
$wsInclude.Range("B4").Value = "http://standardterms.edqm.eu"
